$wb = $excel.ActiveWorkbook

# Update the "thermal" sheet (max_cap for t1 changes from 1000 to 500)
$thermal = $wb.Worksheets.Item("thermal")
$thermal.Range("B2").Value = 500

# Update the "cap_factors" sheet (cap_factor for w1 changes from 0.3 to 0.160507)
$capFactors = $wb.Worksheets.Item("cap_factors")
$capFactors.Range("C2").Value = 0.160507

# Update the "demand" sheet (demand changes from 200 to 246.57255)
$demand = $wb.Worksheets.Item("demand")
$demand.Range("B2").Value = 246.57255

# Select B3 on the thermal sheet and make it the active sheet/tab
$thermal.Activate()
$thermal.Range("B3").Select()

$wb.Save()
